$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is a cell reference + its new display text. A leading apostrophe
# forces Excel to store the value as literal text (matching the source
# workbook's original inline-string cells) instead of inferring a number or
# percentage. Resetting the Style back to "Normal" afterwards drops the
# quote-prefix formatting flag that assigning a leading apostrophe adds, so
# the cell's style index is left exactly as it was before the edit.
$updates = (
    @{ Ref = "D2"; NewValue = "302.99" },
    @{ Ref = "E2"; NewValue = "1.36%" },
    @{ Ref = "D3"; NewValue = "32.53" },
    @{ Ref = "E3"; NewValue = "3.86%" },
    @{ Ref = "D4"; NewValue = "4.932" },
    @{ Ref = "E4"; NewValue = "-3.38%" },
    @{ Ref = "D5"; NewValue = "0.07832" },
    @{ Ref = "E5"; NewValue = "-1.50%" },
    @{ Ref = "D6"; NewValue = "2.040" },
    @{ Ref = "E6"; NewValue = "-10.58%" },
    @{ Ref = "D7"; NewValue = "7.834" },
    @{ Ref = "E7"; NewValue = "0.79%" },
    @{ Ref = "D8"; NewValue = "3.807" },
    @{ Ref = "E8"; NewValue = "-1.41%" },
    @{ Ref = "D9"; NewValue = "0.9210" },
    @{ Ref = "E9"; NewValue = "-0.31%" },
    @{ Ref = "D10"; NewValue = "0.1754" },
    @{ Ref = "E10"; NewValue = "1.16%" },
    @{ Ref = "D11"; NewValue = "0.07887" },
    @{ Ref = "E11"; NewValue = "4.64%" },
    @{ Ref = "D12"; NewValue = "0.08601" },
    @{ Ref = "E12"; NewValue = "-8.99%" },
    @{ Ref = "D13"; NewValue = "0.03138" },
    @{ Ref = "E13"; NewValue = "3.79%" },
    @{ Ref = "D14"; NewValue = "0.1005" },
    @{ Ref = "E14"; NewValue = "0.15%" },
    @{ Ref = "D15"; NewValue = "0.001518" },
    @{ Ref = "E15"; NewValue = "0.66%" },
    @{ Ref = "D16"; NewValue = "0.005902" },
    @{ Ref = "E16"; NewValue = "0.68%" },
    @{ Ref = "D17"; NewValue = "3.471" },
    @{ Ref = "E17"; NewValue = "-0.47%" },
    @{ Ref = "D18"; NewValue = "2.158" },
    @{ Ref = "E18"; NewValue = "-4.82%" },
    @{ Ref = "D19"; NewValue = "0.3308" },
    @{ Ref = "E19"; NewValue = "1.12%" },
    @{ Ref = "D20"; NewValue = "0.1290" },
    @{ Ref = "E20"; NewValue = "-3.27%" },
    @{ Ref = "D21"; NewValue = "4.318" },
    @{ Ref = "E21"; NewValue = "10.60%" },
    @{ Ref = "E22"; NewValue = "17.19%" },
    @{ Ref = "D23"; NewValue = "0.04570" },
    @{ Ref = "E23"; NewValue = "-0.94%" },
    @{ Ref = "E24"; NewValue = "-1.86%" },
    @{ Ref = "E25"; NewValue = "-0.82%" },
    @{ Ref = "E26"; NewValue = "4.20%" },
    @{ Ref = "D39"; NewValue = "0.01741" },
    @{ Ref = "E39"; NewValue = "-1.26%" },
    @{ Ref = "D40"; NewValue = "0.04784" },
    @{ Ref = "E40"; NewValue = "3.86%" },
    @{ Ref = "D41"; NewValue = "0.007479" },
    @{ Ref = "E41"; NewValue = "7.43%" },
    @{ Ref = "D42"; NewValue = "0.1360" },
    @{ Ref = "E42"; NewValue = "0.00%" },
    @{ Ref = "E43"; NewValue = "7.80%" },
    @{ Ref = "E44"; NewValue = "9.53%" },
    @{ Ref = "D45"; NewValue = "0.00006247" },
    @{ Ref = "E45"; NewValue = "-1.02%" },
    @{ Ref = "E46"; NewValue = "0.10%" },
    @{ Ref = "E47"; NewValue = "-61.09%" },
    @{ Ref = "E48"; NewValue = "9.94%" },
    @{ Ref = "E49"; NewValue = "0.10%" },
    @{ Ref = "E50"; NewValue = "0.10%" }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Ref)
    $cell.Value = "'" + $u.NewValue
    $cell.Style = "Normal"
}
